$p = $ppt.ActivePresentation
$p.Slides.Item(6).Delete()
$p.Slides.Item(5).Delete()
